# Apply weekly-bar updates to the BANDHANBNK.NS 1wk sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up existing rows -------------------------------------------------

# Row 325: isPivot (O) flips from 0 -> 2
$ws.Cells.Item(325, 15).Value = 2

# Rows 327-328: "backup" column (R) goes from blank/inline-string to a
# real numeric 0 now that processing has completed for them.
$ws.Cells.Item(327, 18).Value = 0
$ws.Cells.Item(328, 18).Value = 0

# --- Append the newly scraped weekly bars (rows 329-337) -----------------
# Columns: A Datetime, B Open, C High, D Low, E Close, F Adj Close,
#          G Volume, H Year, I Month, J Day, K Hour, L Minute, M Second,
#          N Week, O isPivot, P two_line_structure, Q detect_structure
# Column R (backup) is intentionally left unset for these new rows, matching
# the "not yet processed" state the earlier rows had before this commit.

$newRows = @(
    @{ Row=329; A=45474; B=204;              C=214.6000061035156; D=200.8999938964844; E=204.3300018310547; F=202.7745819091797; G=78655273;  H=2024; I=7; J=1;  N=27; O=0; P=0; Q=0 },
    @{ Row=330; A=45481; B=202.4600067138672; C=205.6000061035156; D=190.0500030517578; E=192.6000061035156; F=191.1338806152344; G=61127651;  H=2024; I=7; J=8;  N=28; O=0; P=0; Q=0 },
    @{ Row=331; A=45488; B=193.3899993896484; C=202.9499969482422; D=191;               E=192.4499969482422; F=190.9850158691406; G=52910353;  H=2024; I=7; J=15; N=29; O=0; P=1; Q=1 },
    @{ Row=332; A=45495; B=191.0299987792969; C=199.4900054931641; D=184.75;             E=192.5;              F=191.0346374511719; G=80388783;  H=2024; I=7; J=22; N=30; O=2; P=0; Q=0 },
    @{ Row=333; A=45502; B=203;               C=222.3099975585938; D=200.1100006103516; E=212.5299987792969; F=210.9121551513672; G=210893028; H=2024; I=7; J=29; N=31; O=0; P=0; Q=0 },
    @{ Row=334; A=45509; B=206.1999969482422; C=209.1999969482422; D=198.1999969482422; E=199.5;              F=197.9813537597656; G=50005064;  H=2024; I=8; J=5;  N=32; O=0; P=0; Q=0 },
    @{ Row=335; A=45516; B=199.1000061035156; C=200;                D=189.5;             E=191.5899963378906; F=190.1315612792969; G=41758804;  H=2024; I=8; J=12; N=33; O=0; P=0; Q=0 },
    @{ Row=336; A=45523; B=193.4499969482422; C=211.0700073242188; D=191.5;              E=201.7200012207031; F=201.7200012207031; G=57723757;  H=2024; I=8; J=19; N=34; O=0; P=0; Q=0 },
    @{ Row=337; A=45530; B=202;               C=203;                D=191.6399993896484; E=200.3000030517578; F=200.3000030517578; G=159013170; H=2024; I=8; J=26; N=35; O=0; P=0; Q=0 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $aCell = $ws.Cells.Item($row, 1)
    $aCell.Value = $r.A
    $aCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = 0
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    # Column R ("backup") deliberately left blank/unset for new rows.
}
